# Remove the override of the attribute ContourVisibility at the Image level.
# (It will be replaced later with an 'Extra Tools' option to toggle the Fill
# option for contours.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shorten the description of the Session-level ContourVisibility attribute
#    now that the Image-level override (and its "See details below" reference)
#    is going away.
$ws.Range("I9").Value = "This defines how contours  (label maps or segmentations) are displayed for the entire quiz - as outline only or as a solid contour."

# The wrapped text in I9 now spans one fewer line, so shrink the row height
# to match (4 lines -> 3 lines at the sheet's default 14.4pt line height).
$ws.Rows.Item(9).RowHeight = 43.2

# 2. Delete the entire "ContourVisibility" row that described the Image-level
#    override (old row 31 - Element=Image, Attribute=ContourVisibility). This
#    removes the attribute name, its options text, and its long note, and
#    shifts every following row up by one.
$ws.Rows.Item(31).Delete()

# 3. Update the active selection/scroll position on the sheet to reflect
#    where the edit was made.
$ws.Range("H10:H11").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 3
